$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $newText
}

# Columns: 1=Variable, 2=N, 3=M, 4=SD, 5=Median, 6=Min, 7=Max

# Row 2: o_llm
Set-CellText $t 2 2 "250"
Set-CellText $t 2 3 "4.78"
Set-CellText $t 2 4 "1.26"
Set-CellText $t 2 5 "4.75"
Set-CellText $t 2 6 "2.42"

# Row 3: c_llm
Set-CellText $t 3 2 "250"
Set-CellText $t 3 3 "4.37"
Set-CellText $t 3 4 "1.18"
Set-CellText $t 3 5 "4.25"
Set-CellText $t 3 7 "7.42"

# Row 4: e_llm
Set-CellText $t 4 2 "250"
Set-CellText $t 4 3 "4.98"
Set-CellText $t 4 4 "1.28"

# Row 5: a_llm
Set-CellText $t 5 2 "250"
Set-CellText $t 5 3 "4.75"
Set-CellText $t 5 4 "1.21"
Set-CellText $t 5 5 "4.75"
Set-CellText $t 5 7 "7.75"

# Row 6: n_llm
Set-CellText $t 6 2 "250"
Set-CellText $t 6 3 "6.40"
Set-CellText $t 6 4 "1.16"
Set-CellText $t 6 6 "2.89"
